$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AI2").Value = 0.88551752005580275
$ws.Range("D3").Value = 0.98784912949161852
$ws.Range("AE3").Value = 0.75811010321917394
$ws.Range("BK3").Value = 0.86754256357487858
$ws.Range("B4").Value = 0.98203063314271222
$ws.Range("BO4").Value = 0.82321507033974972
$ws.Range("BF5").Value = 0.76949212227934294
$ws.Range("AW7").Value = 0.66977369983606649
$ws.Range("AB8").Value = 0.77321143324323682
$ws.Range("G9").Value = 0.92732867960413135
$ws.Range("K9").Value = 0.69674318620072451
$ws.Range("AT9").Value = 0.84064531461389003
$ws.Range("Z10").Value = 0.81170018729228199
$ws.Range("O12").Value = 0.96365757945035402
$ws.Range("R12").Value = 0.98949138222154875
$ws.Range("AN12").Value = 0.73305258706232101
$ws.Range("AZ12").Value = 0.80829068326663744
$ws.Range("BC12").Value = 0.88659861578590138
$ws.Range("K13").Value = 0.89523967034578356
$ws.Range("AZ13").Value = 0.87343803638001616
$ws.Range("B15").Value = 0.88833881781824053
$ws.Range("M15").Value = 0.91492533423496614
$ws.Range("AK15").Value = 0.63148025670726882
$ws.Range("Q16").Value = 0.83929183486949299
$ws.Range("BL16").Value = 0.75272574499692135
$ws.Range("AG17").Value = 0.98173215322926466
$ws.Range("AR17").Value = 0.74855226208937986
$ws.Range("AX17").Value = 0.61712475578201809
$ws.Range("N18").Value = 0.69968637885334206
$ws.Range("S18").Value = 0.90239298459643313
$ws.Range("U19").Value = 0.82264684602023364
$ws.Range("D20").Value = 0.68985349565308995
$ws.Range("R20").Value = 0.91001741079265408
$ws.Range("V20").Value = 0.88090123789898955
$ws.Range("D21").Value = 0.95478961254177275
$ws.Range("F21").Value = 0.70698951386246078
$ws.Range("W21").Value = 0.8939949006864264
$ws.Range("AW21").Value = 0.66121112967056161
$ws.Range("AS22").Value = 0.92665843626149691
$ws.Range("F23").Value = 0.95892930624620942
$ws.Range("N23").Value = 0.81480263304449818
$ws.Range("P23").Value = 0.94508933998498912
$ws.Range("BL23").Value = 0.81503676806368786
$ws.Range("E24").Value = 0.92165048005632966
$ws.Range("O24").Value = 0.9858759565266999
$ws.Range("BH25").Value = 0.87458315470240455
$ws.Range("C26").Value = 0.95417801299537564
$ws.Range("F26").Value = 0.64754035701776624
$ws.Range("V26").Value = 0.98672529177548396
$ws.Range("AA26").Value = 0.83194280897397421
$ws.Range("E27").Value = 0.87526273776567631
$ws.Range("M27").Value = 0.7927174507488054
$ws.Range("U27").Value = 0.803884324992578
$ws.Range("Y27").Value = 0.82590401505009781
$ws.Range("AV27").Value = 0.94220879457661066
$ws.Range("A28").Value = 0.89549528031059378
$ws.Range("P28").Value = 0.61416067464773483
$ws.Range("BO28").Value = 0.84474161602351672
$ws.Range("A29").Value = 0.95370735663616191
$ws.Range("I29").Value = 0.86330938943354107
$ws.Range("W29").Value = 0.80575646658486333
$ws.Range("BM29").Value = 0.90629838128180162
$ws.Range("AD31").Value = 0.64754623679071677
$ws.Range("AK32").Value = 0.72137650852370427
$ws.Range("BB32").Value = 0.92787108054682532
$ws.Range("AI33").Value = 0.93408855746302821
$ws.Range("AV34").Value = 0.60999262549865674
$ws.Range("AD35").Value = 0.96250636705366222
$ws.Range("AN35").Value = 0.78388597761007928
$ws.Range("BC36").Value = 0.8651487595053251
$ws.Range("X37").Value = 0.85468838380307233
$ws.Range("K38").Value = 0.88717469434049256
$ws.Range("AM38").Value = 0.94346187319202035
$ws.Range("AY38").Value = 0.92067458884013786
$ws.Range("K39").Value = 0.79242256914421672
$ws.Range("AW39").Value = 0.95556393735655054
$ws.Range("Y40").Value = 0.67785966827405453
$ws.Range("AM40").Value = 0.83191124056137555
$ws.Range("AX40").Value = 0.8303305727365724
$ws.Range("AM41").Value = 0.99734078586890307
$ws.Range("AN42").Value = 0.97365672183286933
$ws.Range("BE42").Value = 0.95687425137524929
$ws.Range("BL42").Value = 0.85095835958519528
$ws.Range("M43").Value = 0.94634582839124071
$ws.Range("N43").Value = 0.93183659069096003
$ws.Range("AP43").Value = 0.78489874887988065
$ws.Range("D45").Value = 0.93648132889232627
$ws.Range("Q45").Value = 0.73505271921320692
$ws.Range("AR45").Value = 0.749666862512542
$ws.Range("A46").Value = 0.78614529548785272
$ws.Range("Z46").Value = 0.77009475457023502
$ws.Range("AJ46").Value = 0.92464010884009462
$ws.Range("AR46").Value = 0.98029926613787888
$ws.Range("BD46").Value = 0.99442556504542068
$ws.Range("F47").Value = 0.83214371436263324
$ws.Range("K47").Value = 0.83056527251924139
$ws.Range("O48").Value = 0.72869476826049151
$ws.Range("H49").Value = 0.84343370262154216
$ws.Range("AY49").Value = 0.9999886654978174
$ws.Range("AE50").Value = 0.97349536515975466
$ws.Range("BC51").Value = 0.80346147992583539
$ws.Range("AN52").Value = 0.8746940582830447
$ws.Range("AS53").Value = 0.89653096374316488
$ws.Range("BJ53").Value = 0.66381023227740599
$ws.Range("G54").Value = 0.95594605493948437
$ws.Range("AS55").Value = 0.77454458327835074
$ws.Range("BB56").Value = 0.76191691301291109
$ws.Range("BH56").Value = 0.88393868427343114
$ws.Range("BL56").Value = 0.94079403971689546
$ws.Range("AD58").Value = 0.73045837220067944
$ws.Range("AO58").Value = 0.73115546400062348
$ws.Range("C59").Value = 0.9845873651334307
$ws.Range("V59").Value = 0.82616151737317467
$ws.Range("BE59").Value = 0.83282607089477789
$ws.Range("P60").Value = 0.92063908320635957
$ws.Range("AK60").Value = 0.92012051951124085
$ws.Range("AT60").Value = 0.87863399709413637
$ws.Range("AO61").Value = 0.81995951097113295
$ws.Range("BJ61").Value = 0.58893954638086354
$ws.Range("J62").Value = 0.74901745268545983
$ws.Range("AY62").Value = 0.82183673827224379
$ws.Range("BL62").Value = 0.77367896843789175
$ws.Range("AF63").Value = 0.76052966348594775
$ws.Range("AG63").Value = 0.73830158245210664
$ws.Range("BD63").Value = 0.75672582866154969
$ws.Range("BI63").Value = 0.88487308825129007
$ws.Range("AF64").Value = 0.98468020981176108
$ws.Range("BG64").Value = 0.64289274889848991
$ws.Range("AH65").Value = 0.96011713535295096
$ws.Range("O66").Value = 0.77481620576689414
$ws.Range("BM66").Value = 0.74536975499655922
$ws.Range("BF67").Value = 0.7957228018012259
$ws.Range("L68").Value = 0.88324983237038279
$ws.Range("AY68").Value = 0.84165262456046297
$ws.Range("BC68").Value = 0.85877645733104147
$ws.Range("BN68").Value = 0.99374834304126436
